$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = "방향장(direction field)"
$ws.Range("E5").Value = "https://angeloyeo.github.io/2021/04/30/direction_fields.html"

$ws.Range("D9").Value = "MSDA vs. MBA in AI BigData (vs. 일반 MBA)"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/msda-vs-mba-in-ai-bigdata/#utm_source=rss&utm_medium=rss&utm_campaign=msda-vs-mba-in-ai-bigdata"

$ws.Range("D29").Value = "[만화] 인턴일기 13~18"
$ws.Range("E29").Value = "https://blog.promedius.ai/manhwa-inteonilgi-13/"

$ws.Range("D51").Value = "[python] 문자열 내 문자를 모두 대문자 또는 소문자로 바꾸고 싶다면?"
$ws.Range("E51").Value = "https://bskyvision.com/1185"
